$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value2 = 1400.0769  # H18: was 1480
$ws.Cells.Item(18, 10).Value2 = 2000  # J18: was 1999.8
$ws.Cells.Item(18, 12).Value2 = 2000  # L18: was 1999.8
$ws.Cells.Item(18, 14).Value2 = -2568  # N18: was -2567.8
$ws.Cells.Item(76, 8).Value2 = 5679  # H76: was 4239
$ws.Cells.Item(76, 9).Value2 = 6400  # I76: was 4600
$ws.Cells.Item(76, 11).Value2 = 6400  # K76: was 4600
$ws.Cells.Item(76, 13).Value2 = -6085  # M76: was -4285
$ws.Cells.Item(79, 8).Value2 = 5679  # H79: was 4239
$ws.Cells.Item(79, 9).Value2 = 6400  # I79: was 4600
$ws.Cells.Item(79, 11).Value2 = 6400  # K79: was 4600
$ws.Cells.Item(79, 13).Value2 = -5308  # M79: was -3508
$ws.Cells.Item(86, 8).Value2 = 6340  # H86: was 5616.6665
$ws.Cells.Item(86, 9).Value2 = 6975  # I86: was 5980
$ws.Cells.Item(86, 11).Value2 = 6975  # K86: was 5980
$ws.Cells.Item(86, 13).Value2 = -5852  # M86: was -4857
$ws.Cells.Item(89, 8).Value2 = 6340  # H89: was 5616.6665
$ws.Cells.Item(89, 9).Value2 = 6975  # I89: was 5980
$ws.Cells.Item(89, 11).Value2 = 34875  # K89: was 29900
$ws.Cells.Item(89, 13).Value2 = -29259  # M89: was -24284
$ws.Cells.Item(100, 8).Value2 = 2991.2856  # H100: was 2437.5
$ws.Cells.Item(100, 9).Value2 = 2994.75  # I100: was 2437.5
$ws.Cells.Item(100, 10).Value2 = 2986.6667  # J100: was 0
$ws.Cells.Item(100, 11).Value2 = 2994.75  # K100: was 2437.5
$ws.Cells.Item(100, 12).Value2 = 2986.6667  # L100: was 0
$ws.Cells.Item(100, 13).Value2 = -2453.75  # M100: was -1896.5
$ws.Cells.Item(100, 14).Value2 = -4068.6667  # N100: was None

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 5782.1924  # H32: was 6074.9795
$ws.Cells.Item(32, 9).Value2 = 5009.5835  # I32: was 5188.2607
$ws.Cells.Item(32, 10).Value2 = 15053.5  # J32: was 19671.334
$ws.Cells.Item(32, 11).Value2 = 5009.5835  # K32: was 5188.2607
$ws.Cells.Item(32, 12).Value2 = 15053.5  # L32: was 19671.334
$ws.Cells.Item(32, 13).Value2 = -4722.5835  # M32: was -4901.2607
$ws.Cells.Item(32, 14).Value2 = -15627.5  # N32: was -20245.334
$ws.Cells.Item(61, 8).Value2 = 66669252  # H61: was 71431200
$ws.Cells.Item(61, 9).Value2 = 76925600  # I61: was 90911624
$ws.Cells.Item(61, 10).Value2 = 2999.5  # J61: was 2983
$ws.Cells.Item(61, 11).Value2 = 76925600  # K61: was 90911624
$ws.Cells.Item(61, 12).Value2 = 2999.5  # L61: was 2983
$ws.Cells.Item(61, 13).Value2 = -76925388  # M61: was -90911412
$ws.Cells.Item(61, 14).Value2 = -3423.5  # N61: was -3407
$ws.Cells.Item(74, 8).Value2 = 2354.7273  # H74: was 2412.7693
$ws.Cells.Item(74, 9).Value2 = 1878  # I74: was 1937.2
$ws.Cells.Item(74, 10).Value2 = 4500  # J74: was 3998
$ws.Cells.Item(74, 11).Value2 = 1878  # K74: was 1937.2
$ws.Cells.Item(74, 12).Value2 = 4500  # L74: was 3998
$ws.Cells.Item(74, 13).Value2 = -1004  # M74: was -1063.2
$ws.Cells.Item(74, 14).Value2 = -6248  # N74: was -5746
$ws.Cells.Item(77, 8).Value2 = 2354.7273  # H77: was 2412.7693
$ws.Cells.Item(77, 9).Value2 = 1878  # I77: was 1937.2
$ws.Cells.Item(77, 10).Value2 = 4500  # J77: was 3998
$ws.Cells.Item(77, 11).Value2 = 9390  # K77: was 9686
$ws.Cells.Item(77, 12).Value2 = 22500  # L77: was 19990
$ws.Cells.Item(77, 13).Value2 = -5022  # M77: was -5318
$ws.Cells.Item(77, 14).Value2 = -31236  # N77: was -28726
$ws.Cells.Item(97, 8).Value2 = 378.81482  # H97: was 380.7857
$ws.Cells.Item(97, 9).Value2 = 378.81482  # I97: was 380.7857
$ws.Cells.Item(97, 11).Value2 = 378.81482  # K97: was 380.7857
$ws.Cells.Item(97, 13).Value2 = 117.18518  # M97: was 115.2143
$ws.Cells.Item(136, 8).Value2 = 66669252  # H136: was 71431200
$ws.Cells.Item(136, 9).Value2 = 76925600  # I136: was 90911624
$ws.Cells.Item(136, 10).Value2 = 2999.5  # J136: was 2983
$ws.Cells.Item(136, 11).Value2 = 230776800  # K136: was 272734872
$ws.Cells.Item(136, 12).Value2 = 8998.5  # L136: was 8949
$ws.Cells.Item(136, 13).Value2 = -230774250  # M136: was -272732322
$ws.Cells.Item(136, 14).Value2 = -14098.5  # N136: was -14049

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value2 = 25195  # H81: was 58000
$ws.Cells.Item(81, 10).Value2 = 25195  # J81: was 58000
$ws.Cells.Item(81, 12).Value2 = 25195  # L81: was 58000
$ws.Cells.Item(81, 14).Value2 = -27317  # N81: was -60122
$ws.Cells.Item(84, 8).Value2 = 25195  # H84: was 58000
$ws.Cells.Item(84, 10).Value2 = 25195  # J84: was 58000
$ws.Cells.Item(84, 12).Value2 = 75585  # L84: was 174000
$ws.Cells.Item(84, 14).Value2 = -86193  # N84: was -184608
$ws.Cells.Item(94, 8).Value2 = 20834052  # H94: was 25000834
$ws.Cells.Item(94, 9).Value2 = 22727782  # I94: was 27778370
$ws.Cells.Item(94, 11).Value2 = 22727782  # K94: was 27778370
$ws.Cells.Item(94, 13).Value2 = -22727331  # M94: was -27777919
$ws.Cells.Item(105, 8).Value2 = 112210840  # H105: was 91809080
$ws.Cells.Item(105, 9).Value2 = 112210840  # I105: was 91809080
$ws.Cells.Item(105, 11).Value2 = 112210840  # K105: was 91809080
$ws.Cells.Item(105, 13).Value2 = -112209093  # M105: was -91807333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value2 = 345.9  # H22: was 226.28572
$ws.Cells.Item(22, 9).Value2 = 199.66667  # I22: was 171
$ws.Cells.Item(22, 10).Value2 = 408.57144  # J22: was 300
$ws.Cells.Item(22, 11).Value2 = 199.66667  # K22: was 171
$ws.Cells.Item(22, 12).Value2 = 408.57144  # L22: was 300
$ws.Cells.Item(22, 13).Value2 = 150.33333  # M22: was 179
$ws.Cells.Item(22, 14).Value2 = -1108.57144  # N22: was -1000
$ws.Cells.Item(28, 8).Value2 = 20643  # H28: was 17762
$ws.Cells.Item(28, 10).Value2 = 20643  # J28: was 17762
$ws.Cells.Item(28, 12).Value2 = 20643  # L28: was 17762
$ws.Cells.Item(28, 14).Value2 = -21133  # N28: was -18252
$ws.Cells.Item(132, 8).Value2 = 13345.777  # H132: was 19270.334
$ws.Cells.Item(132, 9).Value2 = 17852.334  # I132: was 26406
$ws.Cells.Item(132, 10).Value2 = 4332.6665  # J132: was 4999
$ws.Cells.Item(132, 11).Value2 = 53557.00199999999  # K132: was 79218
$ws.Cells.Item(132, 12).Value2 = 12997.9995  # L132: was 14997
$ws.Cells.Item(132, 13).Value2 = -51027.00199999999  # M132: was -76688
$ws.Cells.Item(132, 14).Value2 = -18057.9995  # N132: was -20057
$ws.Cells.Item(134, 8).Value2 = 33336522  # H134: was 33336588
$ws.Cells.Item(134, 9).Value2 = 3801.2727  # I134: was 4080.2
$ws.Cells.Item(134, 10).Value2 = 125001500  # J134: was 100001600
$ws.Cells.Item(134, 11).Value2 = 11403.8181  # K134: was 12240.6
$ws.Cells.Item(134, 12).Value2 = 375004500  # L134: was 300004800
$ws.Cells.Item(134, 13).Value2 = -8868.8181  # M134: was -9705.599999999999
$ws.Cells.Item(134, 14).Value2 = -375009570  # N134: was -300009870

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(136, 8).Value2 = 1881.6666  # H136: was 1663.6923
$ws.Cells.Item(136, 9).Value2 = 1306  # I136: was 1252.9
$ws.Cells.Item(136, 11).Value2 = 3918  # K136: was 3758.7
$ws.Cells.Item(136, 13).Value2 = 1182  # M136: was 1341.3
$ws.Cells.Item(137, 8).Value2 = 32616728  # H137: was 30007816
$ws.Cells.Item(137, 9).Value2 = 150002400  # I137: was 150003200
$ws.Cells.Item(137, 10).Value2 = 9596.111  # J137: was 8971.2
$ws.Cells.Item(137, 11).Value2 = 450007200  # K137: was 450009600
$ws.Cells.Item(137, 12).Value2 = 28788.333  # L137: was 26913.6
$ws.Cells.Item(137, 13).Value2 = -450002100  # M137: was -450004500
$ws.Cells.Item(137, 14).Value2 = -38988.333  # N137: was -37113.60000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value2 = 893  # H97: was 943
$ws.Cells.Item(97, 9).Value2 = 909.4  # I97: was 943
$ws.Cells.Item(97, 10).Value2 = 811  # J97: was 0
$ws.Cells.Item(97, 11).Value2 = 909.4  # K97: was 943
$ws.Cells.Item(97, 12).Value2 = 811  # L97: was 0
$ws.Cells.Item(97, 13).Value2 = -413.4  # M97: was -447
$ws.Cells.Item(97, 14).Value2 = -1803  # N97: was None
$ws.Cells.Item(102, 8).Value2 = 1150.4  # H102: was 1235.8462
$ws.Cells.Item(102, 9).Value2 = 844  # I102: was 893.8
$ws.Cells.Item(102, 11).Value2 = 844  # K102: was 893.8
$ws.Cells.Item(102, 13).Value2 = 778  # M102: was 728.2
$ws.Cells.Item(122, 8).Value2 = 1605.6875  # H122: was 1655.6875
$ws.Cells.Item(122, 9).Value2 = 1336.8334  # I122: was 1403.4166
$ws.Cells.Item(122, 10).Value2 = 2412.25  # J122: was 2412.5
$ws.Cells.Item(122, 11).Value2 = 4010.5002  # K122: was 4210.2498
$ws.Cells.Item(122, 12).Value2 = 7236.75  # L122: was 7237.5
$ws.Cells.Item(122, 13).Value2 = -1560.5002  # M122: was -1760.2498
$ws.Cells.Item(122, 14).Value2 = -12136.75  # N122: was -12137.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(127, 8).Value2 = 37107.5  # H127: was 0
$ws.Cells.Item(127, 10).Value2 = 37107.5  # J127: was 0
$ws.Cells.Item(127, 12).Value2 = 37107.5  # L127: was 0
$ws.Cells.Item(127, 14).Value2 = -47027.5  # N127: was None
$ws.Cells.Item(132, 8).Value2 = 126722.89  # H132: was 224401.2
$ws.Cells.Item(132, 9).Value2 = 37834.332  # I132: was 55001.5
$ws.Cells.Item(132, 10).Value2 = 171167.17  # J132: was 337334.34
$ws.Cells.Item(132, 11).Value2 = 113502.996  # K132: was 165004.5
$ws.Cells.Item(132, 12).Value2 = 513501.51  # L132: was 1012003.02
$ws.Cells.Item(132, 13).Value2 = -110972.996  # M132: was -162474.5
$ws.Cells.Item(132, 14).Value2 = -518561.51  # N132: was -1017063.02

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value2 = 9262754  # H122: was 9617218
$ws.Cells.Item(122, 9).Value2 = 9262754  # I122: was 10418519
$ws.Cells.Item(122, 10).Value2 = 0  # J122: was 1599.5
$ws.Cells.Item(122, 11).Value2 = 27788262  # K122: was 31255557
$ws.Cells.Item(122, 12).Value2 = 0  # L122: was 4798.5
$ws.Cells.Item(122, 13).Value2 = -27785812  # M122: was -31253107
$ws.Cells.Item(122, 14).ClearContents()  # N122: was -9698.5
$ws.Cells.Item(132, 8).Value2 = 2770.3333  # H132: was 3204.6843
$ws.Cells.Item(132, 9).Value2 = 1233.3334  # I132: was 1581.2
$ws.Cells.Item(132, 10).Value2 = 3692.5334  # J132: was 3784.5
$ws.Cells.Item(132, 11).Value2 = 3700.0002  # K132: was 4743.6
$ws.Cells.Item(132, 12).Value2 = 11077.6002  # L132: was 11353.5
$ws.Cells.Item(132, 13).Value2 = -1170.0002  # M132: was -2213.6
$ws.Cells.Item(132, 14).Value2 = -16137.6002  # N132: was -16413.5
